$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09559038492903756
$ws.Range("C2").Value = 0.9986789149978104
$ws.Range("D2").Value = 0.2511415180113421
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=5, n_estimators=150))])"
$ws.Range("G2").Value = 0.1245324579833929
$ws.Range("H2").Value = 0.992
